# load character data form table
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns (move_speed, movement) before the old F column
# (prefab). This shifts old F,G,H -> H,I,J.
$ws.Range("F1:G1").EntireColumn.Insert()

# Match the column width of the neighboring "group" column (15) for the
# two freshly inserted columns so F:H all render at width 15.
$ws.Columns.Item(6).ColumnWidth = 14.29
$ws.Columns.Item(7).ColumnWidth = 14.29

# Header / type / label rows for the two new columns.
$ws.Range("F1").Value = "move_speed"
$ws.Range("F2").Value = "float"
$ws.Range("F3").Value = "基础移动速度"

$ws.Range("G1").Value = "movement"
$ws.Range("G2").Value = "string"
$ws.Range("G3").Value = "移动方式"

# Data rows for the two new columns.
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = "StraightChase"

$ws.Range("F5").Value = 1
$ws.Range("G5").Value = "StraightChase"

# Update the prefab asset key (now column J) to point at the slime actor.
$ws.Range("J4").Value = "Level:Characters:ActorSlime"
$ws.Range("J5").Value = "Level:Characters:ActorSlime"

# Match the author's final selection (shifted one column right, same row).
$ws.Range("I12").Select()
